$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the bottom of the data (31:32), shifting
# down from the row above so the new rows inherit the existing column
# formatting (email column alignment style, is_active boolean style).
$ws.Range("A31:K32").Insert(-4121)

# Row 32 holds "John Doe" - its strings are entered first so they land
# at the lower shared-string indices.
$ws.Range("A32").Value = 110031
$ws.Range("B32").Value = 9317596767
$ws.Range("C32").Value = "John Doe"
$ws.Range("D32").Value = "john.doe@xyz.com"
$ws.Range("E32").Value = 818876431
$ws.Range("F32").Value = "ACT"
$ws.Range("G32").Value = "eng"
$ws.Range("H32").Value = "PWD"
$ws.Range("I32").Value = $true
$ws.Range("J32").Value = "superadmin"
$ws.Range("K32").Value = "now()"

# Row 31 holds "Jane Smith" - entered after, so its strings take the
# next shared-string slots.
$ws.Range("A31").Value = 110030
$ws.Range("B31").Value = 9317596768
$ws.Range("C31").Value = "Jane Smith"
$ws.Range("D31").Value = "jane.smith@xyz.com"
$ws.Range("E31").Value = 818876432
$ws.Range("F31").Value = "ACT"
$ws.Range("G31").Value = "eng"
$ws.Range("H31").Value = "PWD"
$ws.Range("I31").Value = $true
$ws.Range("J31").Value = "superadmin"
$ws.Range("K31").Value = "now()"

# Match the author's final selection/viewport.
$ws.Range("E28").Select() | Out-Null
